$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This workbook lists bird/lichen observation records (rows 2-13) that share
# the same location cluster. The records were re-synced from source and the
# per-row field values shifted to a new row ordering. Apply the new values
# (Id, Taxonsorteringsordning, TaxonId, Artnamn, Vetenskapligt namn, Auktor,
# Ost, Nord, and the occasional Aktivitet/Publik kommentar note) row by row.

# Row 2
$ws.Range("A2").Value = 106527621
$ws.Range("B2").Value = 56395
$ws.Range("E2").Value = 100109
$ws.Range("F2").Value = "Tretåig hackspett"
$ws.Range("G2").Value = "Picoides tridactylus"
$ws.Range("H2").Value = "(Linnaeus, 1758)"
$ws.Range("Q2").Value = 422551.809521322
$ws.Range("R2").Value = 7019429.901349179
$ws.Range("AC2").Value = "ringhack"

# Row 3
$ws.Range("A3").Value = 106527625
$ws.Range("M3").Value = "färska spår"
$ws.Range("Q3").Value = 422594.7694478295
$ws.Range("R3").Value = 7019417.617479109

# Row 4
$ws.Range("A4").Value = 106527619
$ws.Range("Q4").Value = 422527.6223480041
$ws.Range("R4").Value = 7019435.883743418
$ws.Range("M4").ClearContents()

# Row 5
$ws.Range("A5").Value = 106527614
$ws.Range("Q5").Value = 422426.2143204252
$ws.Range("R5").Value = 7019621.460759236

# Row 6
$ws.Range("A6").Value = 106527638
$ws.Range("B6").Value = 77506
$ws.Range("E6").Value = 6425
$ws.Range("F6").Value = "Garnlav"
$ws.Range("G6").Value = "Alectoria sarmentosa"
$ws.Range("H6").Value = "(Ach.) Ach."
$ws.Range("Q6").Value = 422566.9857729367
$ws.Range("R6").Value = 7019405.237422666
$ws.Range("AC6").ClearContents()

# Row 7
$ws.Range("A7").Value = 106527636
$ws.Range("Q7").Value = 422495.1119518331
$ws.Range("R7").Value = 7019470.414734876

# Row 8
$ws.Range("A8").Value = 106527623
$ws.Range("B8").Value = 56395
$ws.Range("E8").Value = 100109
$ws.Range("F8").Value = "Tretåig hackspett"
$ws.Range("G8").Value = "Picoides tridactylus"
$ws.Range("H8").Value = "(Linnaeus, 1758)"
$ws.Range("Q8").Value = 422587.6292912964
$ws.Range("R8").Value = 7019420.489218531
$ws.Range("AC8").Value = "ringhack"

# Row 9
$ws.Range("A9").Value = 106527634
$ws.Range("B9").Value = 77506
$ws.Range("E9").Value = 6425
$ws.Range("F9").Value = "Garnlav"
$ws.Range("G9").Value = "Alectoria sarmentosa"
$ws.Range("H9").Value = "(Ach.) Ach."
$ws.Range("Q9").Value = 422420.8103881205
$ws.Range("R9").Value = 7019509.553637378
$ws.Range("AC9").ClearContents()

# Row 10
$ws.Range("A10").Value = 106527630
$ws.Range("B10").Value = 77506
$ws.Range("E10").Value = 6425
$ws.Range("F10").Value = "Garnlav"
$ws.Range("G10").Value = "Alectoria sarmentosa"
$ws.Range("H10").Value = "(Ach.) Ach."
$ws.Range("Q10").Value = 421807.3410346356
$ws.Range("R10").Value = 7019388.525846463
$ws.Range("AC10").ClearContents()

# Row 11
$ws.Range("A11").Value = 106527617
$ws.Range("B11").Value = 56395
$ws.Range("E11").Value = 100109
$ws.Range("F11").Value = "Tretåig hackspett"
$ws.Range("G11").Value = "Picoides tridactylus"
$ws.Range("H11").Value = "(Linnaeus, 1758)"
$ws.Range("Q11").Value = 421790.6172308734
$ws.Range("R11").Value = 7019460.47857395
$ws.Range("AC11").Value = "ringhack"

# Row 12
$ws.Range("A12").Value = 106527632
$ws.Range("B12").Value = 77506
$ws.Range("E12").Value = 6425
$ws.Range("F12").Value = "Garnlav"
$ws.Range("G12").Value = "Alectoria sarmentosa"
$ws.Range("H12").Value = "(Ach.) Ach."
$ws.Range("Q12").Value = 421997.6268703607
$ws.Range("R12").Value = 7019616.983763674
$ws.Range("AC12").ClearContents()

# Row 13
$ws.Range("A13").Value = 106527615
$ws.Range("B13").Value = 56395
$ws.Range("E13").Value = 100109
$ws.Range("F13").Value = "Tretåig hackspett"
$ws.Range("G13").Value = "Picoides tridactylus"
$ws.Range("H13").Value = "(Linnaeus, 1758)"
$ws.Range("Q13").Value = 421879.6144842675
$ws.Range("R13").Value = 7019284.621970991
$ws.Range("AC13").Value = "ringhack"
